$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$layout = $master.CustomLayouts.Item(2)
$title = $layout.Shapes.Item(1)
$title.TextFrame.TextRange.Font.Size = 36
